# Daily attendance processing - 2026-01-05 19:57:17
# Re-sort the comma-separated "Recorded By" list (column G) for every data
# row into ordinal (ASCII, case-sensitive) order.

function Sort-Ordinal($arr) {
    $n = $arr.Count
    for ($i = 0; $i -lt $n; $i++) {
        for ($j = 0; $j -lt ($n - $i - 1); $j++) {
            if ($arr[$j].CompareTo($arr[$j + 1]) -gt 0) {
                $tmp = $arr[$j]
                $arr[$j] = $arr[$j + 1]
                $arr[$j + 1] = $tmp
            }
        }
    }
    return $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ur = $ws.UsedRange
$lastRow = $ur.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Value2

    if ($current -eq $null) {
        continue
    }

    $text = [string]$current
    if ($text.IndexOf(",") -lt 0) {
        continue
    }

    $rawParts = $text.Split(",")
    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    $sorted = Sort-Ordinal $parts
    $newText = $sorted -join ", "

    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}
